$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Folder Inventory": the folder
# "Build-Custom-Knowledge-RAG-App-With-Azure-AI-Foundry" was touched again
# (new "Last Updated" timestamp), so after the nightly re-sort (descending by
# Last Updated) it jumps from row 9 up to row 2, pushing the previous rows
# 2-8 down to rows 3-9. Rows 10+ are unaffected.
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Folder Inventory")

$zeroWidthSpace = [char]0x200B

$names = @(
    "Build-Custom-Knowledge-RAG-App-With-Azure-AI-Foundry",
    "Automated Machine Learning Using AML",
    "Create and Publish PowerBI Dashboards & Reports",
    "Azure Virtual Machine And Compute",
    ("Work with Data Lake and Data Factory Pipelines in Microsoft Fabric" + $zeroWidthSpace),
    "Get Started with Microsoft Fabric with Its Lakehouses",
    "Build A Fabric Real-Time Intelligence Solution in a Day",
    "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals"
)

$dates = @(
    "2025-06-13 13:01:36 +0530",
    "2025-06-12 21:50:14 +0530",
    "2025-06-12 20:05:46 +0530",
    "2025-06-12 17:37:08 +0530",
    "2025-06-12 17:26:19 +0530",
    "2025-06-12 16:16:30 +0530",
    "2025-06-12 15:59:35 +0530",
    "2025-06-12 15:19:27 +0530"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $wsInv.Range("A$row").Value = $names[$i]
    $wsInv.Range("B$row").Value = $names[$i]
    $wsInv.Range("C$row").Value = $dates[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Metadata": refresh the generation timestamp and bump the workflow
# run counter.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2025-06-13 07:31:54 UTC"
# "Workflow Run" is stored as text in the workbook (not a number), so force
# a text number format before writing the digit-only string - otherwise
# Excel's COM layer auto-coerces "2" into a numeric value.
$wsMeta.Range("B5").NumberFormat = "@"
$wsMeta.Range("B5").Value = "2"

# ---------------------------------------------------------------------------
# Sheet "Summary": the most recent update now matches the refreshed folder.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-13 13:01:36 +0530"
